$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Pivovar"
$ws.Range("C9").Formula = '="1"'
$ws.Range("C9").Copy()
$ws.Range("C9").PasteSpecial(-4163)
$ws.Range("D9").Value = "marcik@gmail.com"
$ws.Range("E9").Value = "Bruhake"
$ws.Range("F9").Value = "Dance"
$ws.Range("G9").Value = "customer"

$ws.Range("C7").Select()
